$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.392.41'
$ws.Range("E2").Value = '  -0.75%  '
$ws.Range("D3").Value = '1.638.59'
$ws.Range("E3").Value = '  -1.69%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '211.31'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.67%  '
$ws.Range("E6").Value = '  +3.77%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '23.08'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.41%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.256'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.27%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0609'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.12%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0889'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.14%  '
$ws.Range("D12").Value = '1.869.33'
$ws.Range("E12").Value = '  -1.70%  '
$ws.Range("D13").Value = '1.627.41'
$ws.Range("E13").Value = '  -2.27%  '
$ws.Range("E14").Value = '  -2.71%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.558'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.68%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.28'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -2.87%  '
$ws.Range("D17").Value = '27.359.53'
$ws.Range("E17").Value = '  -0.89%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '229.30'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -5.81%  '
$ws.Range("D19").Value = '0.0₃0719'
$ws.Range("E19").Value = '  -1.26%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.57'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.28%  '
$ws.Range("E21").Value = '  +0.08%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.32'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -3.53%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.47'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.13%  '
$ws.Range("E24").Value = '  -0.73%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '148.11'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.94%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.96'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -3.20%  '
$ws.Range("E27").Value = '  +1.36%  '
$ws.Range("E28").Value = '  -0.06%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.54'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -5.54%  '
$ws.Range("E30").Value = '  -4.19%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0484'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -3.28%  '
$ws.Range("E32").Value = '  -2.27%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.12'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.17%  '
$ws.Range("D34").Value = '1.407.09'
$ws.Range("E34").Value = '  -4.07%  '
$ws.Range("E35").Value = '  +0.45%  '
$ws.Range("E36").Value = '  -0.13%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.564'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.58%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.880'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -5.39%  '
$ws.Range("E39").Value = '  -3.56%  '
$ws.Range("E40").Value = '  +1.00%  '
$ws.Range("E41").Value = '  -0.01%  '
$ws.Range("E42").Value = '  -1.59%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.50'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.76%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.24'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.55%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.791'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.46%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '64.47'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -7.11%  '
$ws.Range("D47").Value = '1.779.51'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.64'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -4.60%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '87.32'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.22%  '
$ws.Range("D50").Value = '0.0₆0104'
$ws.Range("E50").Value = '  -3.01%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0989'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.47%  '
